$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("24-1-2022")

$ws.Range("C14").Value = 'having sleep'
$ws.Range("C22").Value = 'now  working on that and also finding how to fix pdf like full fit content so come to concept ppi I increase that also but can''t change anything now I am think increase page height width and reset it again'
$ws.Range("C3").Value = 'so I started a resume project  with new folder new idea in that I one side of page I will be a form and another side has a template and I visits  bunch of sites and see this pattern and observe that I I come to know that its like react project get value from form and render to page of resume'
$ws.Range("C6").Value = 'then get one library which print my page and come to know how to use that library and its work like its read full html CSS and then make pdf like if I added input field so  it not considered the value of that input so I have make this plan'
$ws.Range("C8").Value = 'so I started find a new template in GitHub just for test and print it and watch I can hide any thing or not and do many this like margin and started new template to give normal style '
$ws.Range("C10").Value = 'and then make very simple form to test is value is going and it is replaced and it is replacing make some function on that and printing that'
$ws.Range("C16").Value = 'and then just add new experience and it is working and because of some work I have to go outside'
$ws.Range("C18").Value = 'and there is one problem that there is spacing in the pdf so read its document and find how to zoom it of fit the content but not get it that in documentation nor in stack overflow'
$ws.Range("C20").Value = 'so leave that problem for sometimes  and then add experience and I thing that that’s not going to work but its work then visiting website for observe how I can add bullets of experience and education so I come to know that there is one text area in that there is bullets '

# Restore the tab selection / active cell shown in the workbook view.
$ws.Activate()
$ws.Range("C20").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1

